$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Work Report")

# Update the report generation timestamp
$ws.Range("D5").Value = "Report Generated On: 08/26/2025 10:01 AM"

# Update Total Billed Amount
$ws.Range("C8").Value = 631.14

# Clear Scope ID value (G10)
$ws.Range("G10").Value = ""

# Update Pricing for the line item and TOTAL row
$ws.Range("H16").Value = 631.14
$ws.Range("H17").Value = 631.14
